$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the two new Wormmon skill rows (IDs 60 and 61), matching the
# formatting of the existing data rows above them
$ws.Range("A61:E61").Copy()
$ws.Range("A62:E63").PasteSpecial(-4122)

$ws.Range("A62").Value = 60
$ws.Range("B62").Value = "NaturalSpirit"
$ws.Range("C62").Value = "PassiveSkill(Wormmon)"
$ws.Range("D62").Value = 0
$ws.Range("E62").Value = 0

$ws.Range("A63").Value = 61
$ws.Range("B63").Value = "StickyNet"
$ws.Range("C63").Value = "DamageSkill(Wormmon)"
$ws.Range("D63").Value = 5
$ws.Range("E63").Value = 3

# Match the selection state recorded in the saved file
$ws.Range("D49").Select()
